# Applies the cryptos-list refresh described by the authoritative diff:
# per-row Price (D) / Volume(1h) (E) updates, plus a same-block reordering
# of the Chainlink / WrappedliquidstakedEther2.0 rows (12<->13) and the
# InjectiveProtocol / Kaspa rows (40<->41), each carried with its own
# Coin name, Link and Price/Volume values.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @(
    # Row 2
    @{ Cell = 'D2'; Value = '34.815.30'; ForceText = $false }
    @{ Cell = 'E2'; Value = '  -1.42%  '; ForceText = $false }
    # Row 3
    @{ Cell = 'D3'; Value = '1.868.99'; ForceText = $false }
    @{ Cell = 'E3'; Value = '  -2.21%  '; ForceText = $false }
    # Row 4
    @{ Cell = 'E4'; Value = '  -1.02%  '; ForceText = $false }
    # Row 5
    @{ Cell = 'D5'; Value = '243.49'; ForceText = $true }
    @{ Cell = 'E5'; Value = '  -4.05%  '; ForceText = $false }
    # Row 6
    @{ Cell = 'D6'; Value = '0.671'; ForceText = $true }
    @{ Cell = 'E6'; Value = '  -6.34%  '; ForceText = $false }
    # Row 7
    @{ Cell = 'E7'; Value = '  -1.06%  '; ForceText = $false }
    # Row 8
    @{ Cell = 'D8'; Value = '42.09'; ForceText = $true }
    @{ Cell = 'E8'; Value = '  +3.49%  '; ForceText = $false }
    # Row 9
    @{ Cell = 'D9'; Value = '0.340'; ForceText = $true }
    @{ Cell = 'E9'; Value = '  -5.17%  '; ForceText = $false }
    # Row 10
    @{ Cell = 'D10'; Value = '0.0734'; ForceText = $true }
    @{ Cell = 'E10'; Value = '  -1.87%  '; ForceText = $false }
    # Row 11
    @{ Cell = 'D11'; Value = '0.0965'; ForceText = $true }
    @{ Cell = 'E11'; Value = '  -2.67%  '; ForceText = $false }
    # Row 12
    @{ Cell = 'B12'; Value = 'WrappedliquidstakedEther2.0'; ForceText = $false }
    @{ Cell = 'C12'; Value = 'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth'; ForceText = $false }
    @{ Cell = 'D12'; Value = '2.137.88'; ForceText = $false }
    @{ Cell = 'E12'; Value = '  -2.31%  '; ForceText = $false }
    # Row 13
    @{ Cell = 'B13'; Value = 'Chainlink'; ForceText = $false }
    @{ Cell = 'C13'; Value = 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'; ForceText = $false }
    @{ Cell = 'D13'; Value = '12.77'; ForceText = $true }
    @{ Cell = 'E13'; Value = '  +1.39%  '; ForceText = $false }
    # Row 14
    @{ Cell = 'D14'; Value = '0.708'; ForceText = $true }
    @{ Cell = 'E14'; Value = '  -1.33%  '; ForceText = $false }
    # Row 15
    @{ Cell = 'D15'; Value = '1.869.81'; ForceText = $false }
    @{ Cell = 'E15'; Value = '  -2.06%  '; ForceText = $false }
    # Row 16
    @{ Cell = 'D16'; Value = '4.81'; ForceText = $true }
    @{ Cell = 'E16'; Value = '  -1.93%  '; ForceText = $false }
    # Row 17
    @{ Cell = 'D17'; Value = '34.791.13'; ForceText = $false }
    @{ Cell = 'E17'; Value = '  -1.52%  '; ForceText = $false }
    # Row 18
    @{ Cell = 'D18'; Value = '72.09'; ForceText = $true }
    @{ Cell = 'E18'; Value = '  -2.76%  '; ForceText = $false }
    # Row 19
    @{ Cell = 'D19'; Value = '0.0₃0808'; ForceText = $false }
    @{ Cell = 'E19'; Value = '  -4.52%  '; ForceText = $false }
    # Row 20
    @{ Cell = 'D20'; Value = '242.15'; ForceText = $true }
    @{ Cell = 'E20'; Value = '  -0.64%  '; ForceText = $false }
    # Row 21
    @{ Cell = 'D21'; Value = '12.53'; ForceText = $true }
    @{ Cell = 'E21'; Value = '  -3.37%  '; ForceText = $false }
    # Row 22
    @{ Cell = 'D22'; Value = '4.87'; ForceText = $true }
    @{ Cell = 'E22'; Value = '  -3.84%  '; ForceText = $false }
    # Row 23
    @{ Cell = 'E23'; Value = '  -0.94%  '; ForceText = $false }
    # Row 24
    @{ Cell = 'D24'; Value = '2.48'; ForceText = $true }
    @{ Cell = 'E24'; Value = '  +5.30%  '; ForceText = $false }
    # Row 25
    @{ Cell = 'D25'; Value = '2.13'; ForceText = $true }
    @{ Cell = 'E25'; Value = '  -13.10%  '; ForceText = $false }
    # Row 26
    @{ Cell = 'D26'; Value = '162.97'; ForceText = $true }
    @{ Cell = 'E26'; Value = '  -2.29%  '; ForceText = $false }
    # Row 27
    @{ Cell = 'D27'; Value = '8.30'; ForceText = $true }
    @{ Cell = 'E27'; Value = '  -3.51%  '; ForceText = $false }
    # Row 28
    @{ Cell = 'D28'; Value = '18.00'; ForceText = $true }
    @{ Cell = 'E28'; Value = '  -3.85%  '; ForceText = $false }
    # Row 29
    @{ Cell = 'E29'; Value = '  -5.69%  '; ForceText = $false }
    # Row 30
    @{ Cell = 'D30'; Value = '4.128.46'; ForceText = $false }
    @{ Cell = 'E30'; Value = '  +0.05%  '; ForceText = $false }
    # Row 31
    @{ Cell = 'D31'; Value = '1.72'; ForceText = $true }
    @{ Cell = 'E31'; Value = '  +5.60%  '; ForceText = $false }
    # Row 32
    @{ Cell = 'D32'; Value = '4.16'; ForceText = $true }
    @{ Cell = 'E32'; Value = '  -4.30%  '; ForceText = $false }
    # Row 33
    @{ Cell = 'E33'; Value = '  -2.49%  '; ForceText = $false }
    # Row 34
    @{ Cell = 'E34'; Value = '  -1.09%  '; ForceText = $false }
    # Row 35
    @{ Cell = 'E35'; Value = '  -2.27%  '; ForceText = $false }
    # Row 36
    @{ Cell = 'E36'; Value = '  -9.74%  '; ForceText = $false }
    # Row 37
    @{ Cell = 'D37'; Value = '1.93'; ForceText = $true }
    @{ Cell = 'E37'; Value = '  -4.38%  '; ForceText = $false }
    # Row 38
    @{ Cell = 'E38'; Value = '  -25.80%  '; ForceText = $false }
    # Row 39
    @{ Cell = 'D39'; Value = '97.27'; ForceText = $true }
    @{ Cell = 'E39'; Value = '  +0.33%  '; ForceText = $false }
    # Row 40
    @{ Cell = 'B40'; Value = 'Kaspa'; ForceText = $false }
    @{ Cell = 'C40'; Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'; ForceText = $false }
    @{ Cell = 'D40'; Value = '0.0663'; ForceText = $true }
    @{ Cell = 'E40'; Value = '  +2.60%  '; ForceText = $false }
    # Row 41
    @{ Cell = 'B41'; Value = 'InjectiveProtocol'; ForceText = $false }
    @{ Cell = 'C41'; Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'; ForceText = $false }
    @{ Cell = 'D41'; Value = '16.74'; ForceText = $true }
    @{ Cell = 'E41'; Value = '  -2.97%  '; ForceText = $false }
    # Row 42
    @{ Cell = 'E42'; Value = '  -4.25%  '; ForceText = $false }
    # Row 43
    @{ Cell = 'E43'; Value = '  -3.83%  '; ForceText = $false }
    # Row 44
    @{ Cell = 'D44'; Value = '0.0824'; ForceText = $true }
    @{ Cell = 'E44'; Value = '  +11.47%  '; ForceText = $false }
    # Row 45
    @{ Cell = 'D45'; Value = '1.280.26'; ForceText = $false }
    @{ Cell = 'E45'; Value = '  -4.42%  '; ForceText = $false }
    # Row 46
    @{ Cell = 'E46'; Value = '  -5.44%  '; ForceText = $false }
    # Row 47
    @{ Cell = 'E47'; Value = '  -1.35%  '; ForceText = $false }
    # Row 48
    @{ Cell = 'E48'; Value = '  -1.90%  '; ForceText = $false }
    # Row 49
    @{ Cell = 'D49'; Value = '11.80'; ForceText = $true }
    @{ Cell = 'E49'; Value = '  -1.97%  '; ForceText = $false }
    # Row 50
    @{ Cell = 'D50'; Value = '6.25'; ForceText = $true }
    @{ Cell = 'E50'; Value = '  -7.79%  '; ForceText = $false }
    # Row 51
    @{ Cell = 'D51'; Value = '42.48'; ForceText = $true }
    @{ Cell = 'E51'; Value = '  -5.72%  '; ForceText = $false }
)

foreach ($u in $updates) {
    $rng = $ws.Range($u.Cell)
    if ($u.ForceText) {
        # Source cell values that look numeric (e.g. '243.49') must stay
        # plain text, matching the original inline-string-typed column --
        # set the cell's number format to Text before assigning so Excel
        # doesn't silently coerce it to a Number.
        $rng.NumberFormat = '@'
    }
    $rng.Value = $u.Value
}
